$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.422.60'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.099.53'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.34'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.51%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5227'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.56%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4565'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.92%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '56.65'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +15.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08935'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.177'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.87%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.19'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.097.75'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.838'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.051'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '97.37'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001151'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.75%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06635'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.19'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.54%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.307'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.473.64'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.28%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.358'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.336.37'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.19'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.91'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.515'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.28'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.209'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1069'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.659'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.357'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.944'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.30'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.914'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.91%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06858'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2325'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.58%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6871'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.246'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.23%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.27%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6397'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.65%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '14.01'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.53%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.61%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.41%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.06'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +14.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.200'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.58%  '
